$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data describing additional FE data sources, plus a final
# "timeline" style entry (row 8) with only column A populated.
$ws.Range("A3").Value = "Subnational indicators explorer "
$ws.Range("B3").Value = "Compare a local authority and the UK average (median) local authority by different indicators, such as weekly pay and healthy life expectancy."
$ws.Range("C3").Value = "ONS"
$ws.Range("D3").Value = "https://www.ons.gov.uk/peoplepopulationandcommunity/wellbeing/articles/subnationalindicatorsexplorer/2022-01-06"
$ws.Range("E3").Value = "Publicly available"

$ws.Range("A4").Value = "Association of Colleges (AoC)"
$ws.Range("B4").Value = "Lists useful data sources related to colleges"
$ws.Range("C4").Value = "AoC"
$ws.Range("D4").Value = "https://www.aoc.co.uk/research-unit/data-sources"
$ws.Range("E4").Value = "Publicly available"

$ws.Range("A5").Value = "Census"
$ws.Range("B5").Value = "Contains a range of topics including labour market and education by local authority. "
$ws.Range("C5").Value = "ONS"
$ws.Range("D5").Value = "https://census.gov.uk/local-authorities"
$ws.Range("E5").Value = "Publicly available"

$ws.Range("A6").Value = "NOMIS"
$ws.Range("B6").Value = "Query labour market data at national, regional and local levels"
$ws.Range("C6").Value = "ONS"
$ws.Range("D6").Value = "https://www.nomisweb.co.uk/"
$ws.Range("E6").Value = "Publicly available, with more functionality if you create an account. "

$ws.Range("A7").Value = "Explore Education Statistics"
$ws.Range("B7").Value = "A range of educational data published by the Department of Education "
$ws.Range("C7").Value = "DfE"
$ws.Range("D7").Value = "https://explore-education-statistics.service.gov.uk/"
$ws.Range("E7").Value = "Publicly available "

$ws.Range("A8").Value = "Working Futures 20235"

# Column width adjustments to fit the new, wider content. (ColumnWidth is in
# "characters"; Excel stores the OOXML <col width> in a quantized pixel grid,
# so we feed back the target minus the fixed 5px/MDW padding the host adds.)
$ws.Columns.Item(1).ColumnWidth = 12.592447916666666
$ws.Columns.Item(2).ColumnWidth = 39.307291666666664
$ws.Columns.Item(3).ColumnWidth = 5.592447916666667
$ws.Columns.Item(4).ColumnWidth = 11.592447916666666
$ws.Columns.Item(5).ColumnWidth = 19.877604166666668

# Move the active selection the way the author's session left it.
$ws.Range("B14").Select()
